$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.189.96'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.110.44'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.45%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '624.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.88%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.380'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.972'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +24.00%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.106.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.727'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +23.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.191'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.87%  '
$ws.Range("E13").Value = '  +9.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.49'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.42'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.999.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.682.18'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.112.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.99%  '
$ws.Range("E19").Value = '  +15.19%  '
$ws.Range("E20").Value = '  +12.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '436.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.21'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '86.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.273.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.29%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  -4.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +14.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '525.39'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.902'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -16.86%  '
$ws.Range("E34").Value = '  +6.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.144'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '23.70'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.24%  '
$ws.Range("E38").Value = '  +5.37%  '
$ws.Range("E39").Value = '  +5.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0851'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +24.82%  '

$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.20%  '

$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.400'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +12.26%  '

$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.150'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +16.38%  '

$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("E46").Value = '  +8.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '146.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.69%  '
$ws.Range("E49").Value = '  +10.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '167.54'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.42%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.21'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.56%  '

